$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.951.61"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.638.16"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +1.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.85"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "1.865.11"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.25"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "1.620.61"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "25.957.67"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.14"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.25"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.54"
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.903"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").Value = "1.138.72"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.545"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.27"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.41"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").Value = "1.774.73"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("E45").Value = "  +7.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.47"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -1.15%  "
